# live_trading_results.xlsx — Trade #123 closed at 2026-02-16 21:45:28 - leadlag UP +0.000%
#
# This applies:
#   1. Summary sheet: refreshed aggregate stats (OVERALL + leadlag rows) now
#      that two more leadlag trades (#102, #103) have closed.
#   2. leadlag sheet: trades #102 (row 78) and #103 (row 79) flip from
#      OPEN -> CLOSED with their exit data filled in; a brand new trade
#      #123 (row 98) is appended as OPEN.
#   3. All Trades sheet: the two newly-closed trades are appended as rows
#      103 and 104 (this sheet only carries CLOSED trades).
#   4. Comparison sheet: leadlag row refreshed to match the new totals.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    # Writes $Value to $Range as literal text, even when it looks like a
    # number/date/percentage (e.g. "72.8%", "2026-02-16"), and leaves the
    # cell's style untouched (matches the workbook's existing plain
    # inlineStr/shared-string cells, which carry no explicit format).
    param($Range, $Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1. Summary sheet — OVERALL row (row 2) and leadlag row (row 3)
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

$summary.Cells.Item(2, 3).Value = 103
Set-TextValue $summary.Cells.Item(2, 4) "72.8%"
Set-TextValue $summary.Cells.Item(2, 5) "+33.0413%"
Set-TextValue $summary.Cells.Item(2, 6) "+0.3208%"

$summary.Cells.Item(3, 3).Value = 96
Set-TextValue $summary.Cells.Item(3, 4) "54.2%"
Set-TextValue $summary.Cells.Item(3, 5) "+19.2307%"
Set-TextValue $summary.Cells.Item(3, 6) "+0.2003%"

# ---------------------------------------------------------------------
# 2. leadlag sheet — close trade #102 (row 78) and #103 (row 79),
#    append new trade #123 (row 98)
# ---------------------------------------------------------------------
$leadlag = $wb.Worksheets.Item("leadlag")

# Trade #102 -> CLOSED
$leadlag.Cells.Item(78, 7).Value = 69032.605408
$leadlag.Cells.Item(78, 8).Value = "CLOSED"
$leadlag.Cells.Item(78, 9).Value = 0.9068000000000001
$leadlag.Cells.Item(78, 10).Value = 9.07
Set-TextValue $leadlag.Cells.Item(78, 13) "time_exit_5min"
$leadlag.Cells.Item(78, 14).Value = 5

# Trade #103 -> CLOSED
$leadlag.Cells.Item(79, 7).Value = 68676.76274400001
$leadlag.Cells.Item(79, 8).Value = "CLOSED"
$leadlag.Cells.Item(79, 9).Value = 0.3597
$leadlag.Cells.Item(79, 10).Value = 3.6
Set-TextValue $leadlag.Cells.Item(79, 13) "time_exit_5min"
$leadlag.Cells.Item(79, 14).Value = 5

# New trade #123 -> OPEN (row 98)
$leadlag.Cells.Item(98, 1).Value = 123
Set-TextValue $leadlag.Cells.Item(98, 2) "2026-02-16"
Set-TextValue $leadlag.Cells.Item(98, 3) "21:45:28"
$leadlag.Cells.Item(98, 4).Value = "leadlag"
$leadlag.Cells.Item(98, 5).Value = "UP"
$leadlag.Cells.Item(98, 6).Value = 68435.55499999999
$leadlag.Cells.Item(98, 8).Value = "OPEN"
$leadlag.Cells.Item(98, 9).Value = 0
$leadlag.Cells.Item(98, 10).Value = 0
$leadlag.Cells.Item(98, 11).Value = 0.75
$leadlag.Cells.Item(98, 12).Value = "Binance leading with 0.127% move"
$leadlag.Cells.Item(98, 14).Value = 0

# ---------------------------------------------------------------------
# 3. All Trades sheet — append the two newly-closed trades as rows 103/104
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Row 103 <- leadlag trade #102
$allTrades.Cells.Item(103, 1).Value = 102
Set-TextValue $allTrades.Cells.Item(103, 2) "2026-02-16"
Set-TextValue $allTrades.Cells.Item(103, 3) "21:40:21"
$allTrades.Cells.Item(103, 4).Value = "leadlag"
$allTrades.Cells.Item(103, 5).Value = "UP"
$allTrades.Cells.Item(103, 6).Value = 68412.245
$allTrades.Cells.Item(103, 7).Value = 69032.605408
$allTrades.Cells.Item(103, 8).Value = "CLOSED"
$allTrades.Cells.Item(103, 9).Value = 0.9068000000000001
$allTrades.Cells.Item(103, 10).Value = 9.07
$allTrades.Cells.Item(103, 11).Value = 0.75
$allTrades.Cells.Item(103, 12).Value = "Binance leading with 0.146% move"
Set-TextValue $allTrades.Cells.Item(103, 13) "time_exit_5min"
$allTrades.Cells.Item(103, 14).Value = 5

# Row 104 <- leadlag trade #103
$allTrades.Cells.Item(104, 1).Value = 103
Set-TextValue $allTrades.Cells.Item(104, 2) "2026-02-16"
Set-TextValue $allTrades.Cells.Item(104, 3) "21:40:28"
$allTrades.Cells.Item(104, 4).Value = "leadlag"
$allTrades.Cells.Item(104, 5).Value = "UP"
$allTrades.Cells.Item(104, 6).Value = 68430.595
$allTrades.Cells.Item(104, 7).Value = 68676.76274400001
$allTrades.Cells.Item(104, 8).Value = "CLOSED"
$allTrades.Cells.Item(104, 9).Value = 0.3597
$allTrades.Cells.Item(104, 10).Value = 3.6
$allTrades.Cells.Item(104, 11).Value = 0.75
$allTrades.Cells.Item(104, 12).Value = "Coinbase leading with 0.101% move"
Set-TextValue $allTrades.Cells.Item(104, 13) "time_exit_5min"
$allTrades.Cells.Item(104, 14).Value = 5

# ---------------------------------------------------------------------
# 4. Comparison sheet — leadlag row (row 2) refreshed
# ---------------------------------------------------------------------
$comparison = $wb.Worksheets.Item("Comparison")

$comparison.Cells.Item(2, 2).Value = 96
Set-TextValue $comparison.Cells.Item(2, 3) "54.2%"
Set-TextValue $comparison.Cells.Item(2, 4) "3.48"
Set-TextValue $comparison.Cells.Item(2, 5) "+0.5188%"
Set-TextValue $comparison.Cells.Item(2, 7) "1.74"
